$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.04339299999999999
$ws.Range("H2").Value = 0.130179
$ws.Range("I2").Value = 0.0698021577815419
$ws.Range("J2").Value = 0.0698021577815419
$ws.Range("M2").Value = 0.295943
$ws.Range("N2").Value = 0.887829
$ws.Range("O2").Value = 0.0553754985456454
$ws.Range("P2").Value = 0.0553754985456454
$ws.Range("Q2").Value = 0.012841854599
$ws.Range("R2").Value = 0.115576691391
$ws.Range("S2").Value = 0.003865329286714684
$ws.Range("T2").Value = 0.003865329286714684

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.04339299999999999
$ws.Range("H3").Value = 0.130179
$ws.Range("I3").Value = 0.0698021577815419
$ws.Range("J3").Value = 0.0698021577815419
$ws.Range("O3").Value = 0.9446245014543545
$ws.Range("P3").Value = 0.9446245014543546
$ws.Range("Q3").Value = 0.219063138336
$ws.Range("R3").Value = 1.971568245024
$ws.Range("S3").Value = 0.06593682849482721
$ws.Range("T3").Value = 0.06593682849482722

# Row 4
$ws.Range("I4").Value = 0.6764796878879081
$ws.Range("J4").Value = 0.6764796878879081
$ws.Range("M4").Value = 0.295943
$ws.Range("N4").Value = 0.887829
$ws.Range("O4").Value = 0.0553754985456454
$ws.Range("P4").Value = 0.0553754985456454
$ws.Range("Q4").Value = 0.1244553759816667
$ws.Range("R4").Value = 1.120098383835
$ws.Range("S4").Value = 0.03746039997279551
$ws.Range("T4").Value = 0.03746039997279551

# Row 5
$ws.Range("I5").Value = 0.6764796878879081
$ws.Range("J5").Value = 0.6764796878879081
$ws.Range("O5").Value = 0.9446245014543545
$ws.Range("P5").Value = 0.9446245014543546
$ws.Range("S5").Value = 0.6390192879151125
$ws.Range("T5").Value = 0.6390192879151126

# Row 6
$ws.Range("I6").Value = 0.2537181543305499
$ws.Range("J6").Value = 0.2537181543305499
$ws.Range("M6").Value = 0.295943
$ws.Range("N6").Value = 0.887829
$ws.Range("O6").Value = 0.0553754985456454
$ws.Range("P6").Value = 0.0553754985456454
$ws.Range("Q6").Value = 0.04667780697033334
$ws.Range("R6").Value = 0.420100262733
$ws.Range("S6").Value = 0.0140497692861352
$ws.Range("T6").Value = 0.0140497692861352

# Row 7
$ws.Range("I7").Value = 0.2537181543305499
$ws.Range("J7").Value = 0.2537181543305499
$ws.Range("O7").Value = 0.9446245014543545
$ws.Range("P7").Value = 0.9446245014543546
$ws.Range("S7").Value = 0.2396683850444147
$ws.Range("T7").Value = 0.2396683850444148
